$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 38, shifting existing rows 38-66 down to 39-67.
$ws.Rows(38).Insert()

# Populate the newly inserted row 38 with the new data record.
$ws.Cells.Item(38, 1).Value = 10
$ws.Cells.Item(38, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(38, 3).Value = "La Araucanía"
$ws.Cells.Item(38, 4).Value = 44596
$ws.Cells.Item(38, 5).Value = 9
$ws.Cells.Item(38, 6).Value = 100112030
$ws.Cells.Item(38, 7).Value = "Poroto granado"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 40
$ws.Cells.Item(38, 11).Value = 28000
$ws.Cells.Item(38, 12).Value = 28000
$ws.Cells.Item(38, 13).Value = 28000
$ws.Cells.Item(38, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(38, 15).Value = "Región del Maule"
$ws.Cells.Item(38, 16).Value = 1120
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"
